$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# Copy the number formatting (styles) used by the preceding data row (112)
# down onto the 12 new rows (113-124) so the new cells match the existing
# date/time (col A) and fractional-time (col B) formats exactly.
$ws.Range("A112:B112").Copy()
$ws.Range("A113:B124").PasteSpecial(-4122)

# New weigh-in log entries (most-recent-first log, appended at the bottom)
$ws.Cells.Item(113, 1).Value = 44089.36041666667
$ws.Cells.Item(113, 2).Value = 0.36041666666666666
$ws.Cells.Item(113, 3).Value = 72.2

$ws.Cells.Item(114, 1).Value = 44089.359722222223
$ws.Cells.Item(114, 2).Value = 0.35972222222222222
$ws.Cells.Item(114, 3).Value = 72.2

$ws.Cells.Item(115, 1).Value = 44089.336111111108
$ws.Cells.Item(115, 2).Value = 0.33611111111111108
$ws.Cells.Item(115, 3).Value = 72.2

$ws.Cells.Item(116, 1).Value = 44088.879166666666
$ws.Cells.Item(116, 2).Value = 0.87916666666666676
$ws.Cells.Item(116, 3).Value = 72.7

$ws.Cells.Item(117, 1).Value = 44088.493750000001
$ws.Cells.Item(117, 2).Value = 0.49374999999999997
$ws.Cells.Item(117, 3).Value = 72.3

$ws.Cells.Item(118, 1).Value = 44088.338888888888
$ws.Cells.Item(118, 2).Value = 0.33888888888888885
$ws.Cells.Item(118, 3).Value = 72.3

$ws.Cells.Item(119, 1).Value = 44088.338194444441
$ws.Cells.Item(119, 2).Value = 0.33819444444444446
$ws.Cells.Item(119, 3).Value = 72.8

$ws.Cells.Item(120, 1).Value = 44088.308333333334
$ws.Cells.Item(120, 2).Value = 0.30833333333333335
$ws.Cells.Item(120, 3).Value = 72.8

$ws.Cells.Item(121, 1).Value = 44088.281944444447
$ws.Cells.Item(121, 2).Value = 0.28194444444444444
$ws.Cells.Item(121, 3).Value = 73.099999999999994

$ws.Cells.Item(122, 1).Value = 44087.922222222223
$ws.Cells.Item(122, 2).Value = 0.92222222222222217
$ws.Cells.Item(122, 3).Value = 73.7

$ws.Cells.Item(123, 1).Value = 44087.48333333333
$ws.Cells.Item(123, 2).Value = 0.48333333333333334
$ws.Cells.Item(123, 3).Value = 72

$ws.Cells.Item(124, 1).Value = 44087.270833333336
$ws.Cells.Item(124, 2).Value = 0.27083333333333331
$ws.Cells.Item(124, 3).Value = 72.5

# Extend the AM/PM helper formula from column D down through the new rows
$ws.Range("D113:D124").Formula = "=IF(B113<TIME(12,0,0), ""AM"", ""PM"")"

# Update the view: selection moves to the first empty row below the new data
$ws.Activate()
$ws.Range("A125").Select()
